$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-17 (A:T), reflecting the new M2 sending/target cluster
# added to the F2 -> Thbd ligand-receptor pair table.
$data = @(
    @("ECs","F2","Thbd","ECs",3,1,0.6740386666666667,2.022116,0.3058319194585966,0.3058319194585966,3,1,44.32293066666667,132.968792,0.6264651541551962,0.6264651541551961,29.87536908931911,268.878321803872,0.1915930405692093,0.1915930405692093),
    @("ECs","F2","Thbd","FAPs",3,1,0.6740386666666667,2.022116,0.3058319194585966,0.3058319194585966,3,1,17.71526766666667,53.145803,0.2503895325235164,0.2503895325235164,11.94077539768311,107.466978579148,0.07657711134400771,0.07657711134400771),
    @("ECs","F2","Thbd","M2",3,1,0.6740386666666667,2.022116,0.3058319194585966,0.3058319194585966,3,1,4.575739333333334,13.727218,0.06467400065191224,0.06467400065191224,3.084225239254223,27.758027153288,0.01977937375844085,0.01977937375844085),
    @("ECs","F2","Thbd","sCs",3,1,0.6740386666666667,2.022116,0.3058319194585966,0.3058319194585966,3,1,4.136894000000001,12.410682,0.05847131266937523,0.05847131266937522,2.788426515901334,25.095838643112,0.01788239378693878,0.01788239378693878),
    @("FAPs","F2","Thbd","ECs",3,1,0.9731926666666667,2.919578,0.4415672215387696,0.4415672215387696,3,1,44.32293066666667,132.968792,0.6264651541551962,0.6264651541551961,43.13475108997512,388.212759809776,0.276626477511167,0.2766264775111669),
    @("FAPs","F2","Thbd","FAPs",3,1,0.9731926666666667,2.919578,0.4415672215387696,0.4415672215387696,3,1,17.71526766666667,53.145803,0.2503895325235164,0.2503895325235164,17.24036858123711,155.163317231134,0.1105638101788005,0.1105638101788005),
    @("FAPs","F2","Thbd","M2",3,1,0.9731926666666667,2.919578,0.4415672215387696,0.4415672215387696,3,1,4.575739333333334,13.727218,0.06467400065191224,0.06467400065191224,4.453075963778223,40.077683674004,0.02855791877366146,0.02855791877366146),
    @("FAPs","F2","Thbd","sCs",3,1,0.9731926666666667,2.919578,0.4415672215387696,0.4415672215387696,3,1,4.136894000000001,12.410682,0.05847131266937523,0.05847131266937522,4.025994903577335,36.233954132196,0.02581901507514068,0.02581901507514067),
    @("M2","F2","Thbd","ECs",1,0.3333333333333333,0.048595,0.145785,0.02204903496054208,0.02204903496054208,3,1,44.32293066666667,132.968792,0.6264651541551962,0.6264651541551961,2.153872815746667,19.38485534172,0.0138129520855293,0.0138129520855293),
    @("M2","F2","Thbd","FAPs",1,0.3333333333333333,0.048595,0.145785,0.02204903496054208,0.02204903496054208,3,1,17.71526766666667,53.145803,0.2503895325235164,0.2503895325235164,0.8608734322616666,7.747860890355,0.0055208475563648,0.0055208475563648),
    @("M2","F2","Thbd","M2",1,0.3333333333333333,0.048595,0.145785,0.02204903496054208,0.02204903496054208,3,1,4.575739333333334,13.727218,0.06467400065191224,0.06467400065191224,0.2223580529033334,2.00122247613,0.001425999301412134,0.001425999301412134),
    @("M2","F2","Thbd","sCs",1,0.3333333333333333,0.048595,0.145785,0.02204903496054208,0.02204903496054208,3,1,4.136894000000001,12.410682,0.05847131266937523,0.05847131266937522,0.20103236393,1.80929127537,0.001289236017235841,0.001289236017235841),
    @("sCs","F2","Thbd","ECs",3,1,0.508125,1.524375,0.2305518240420917,0.2305518240420917,3,1,44.32293066666667,132.968792,0.6264651541551962,0.6264651541551961,22.521589145,202.694302305,0.1444326839892906,0.1444326839892906),
    @("sCs","F2","Thbd","FAPs",3,1,0.508125,1.524375,0.2305518240420917,0.2305518240420917,3,1,17.71526766666667,53.145803,0.2503895325235164,0.2503895325235164,9.001570383125001,81.01413344812501,0.05772776344434333,0.05772776344434333),
    @("sCs","F2","Thbd","M2",3,1,0.508125,1.524375,0.2305518240420917,0.2305518240420917,3,1,4.575739333333334,13.727218,0.06467400065191224,0.06467400065191224,2.325047548750001,20.92542793875,0.01491070881839779,0.01491070881839779),
    @("sCs","F2","Thbd","sCs",3,1,0.508125,1.524375,0.2305518240420917,0.2305518240420917,3,1,4.136894000000001,12.410682,0.05847131266937523,0.05847131266937522,2.102059263750001,18.91853337375,0.01348066779005992,0.01348066779005992)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
